# Updated slides and exercises
#
# 1) Bump the cached "last modified" date field shown in the footer of all
#    three slide masters from 27.02.2020 to 28.02.2020.
# 2) Add a space before the opening brace in a handful of Java code
#    snippets ("(){" -> "() {", "){" -> ") {", "()){" -> "()) {").

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide-master footer date fields (one shape per master, all identical
#    cached text "27.02.2020").
# ---------------------------------------------------------------------
for ($i = 1; $i -le $p.Designs.Count; $i++) {
    $master = $p.Designs.Item($i).SlideMaster
    for ($j = 1; $j -le $master.Shapes.Count; $j++) {
        $shape = $master.Shapes.Item($j)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq "27.02.2020") {
                $tr.Text = "28.02.2020"
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Java code snippets: insert a space before "{" after a parameter list.
# ---------------------------------------------------------------------

# Slide 9, shape 3: "public EchoFormPage(){" -> "public EchoFormPage() {"
$sh = $p.Slides.Item(9).Shapes.Item(3)
$tr = $sh.TextFrame.TextRange
$f = $tr.Find("(){")
$f.Text = "() {"

# Slide 22, shape 2: three edits inside the same code listing.
$sh = $p.Slides.Item(22).Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# "public ClickCounterPage(){" -> "public ClickCounterPage() {"
$f = $tr.Find("(){")
$f.Text = "() {"

# "link = new Link<Void>(""link-link""){" -> "... ) {"
# (the "(){" occurrence above was already fixed, so the first remaining
# plain "){" match is this one; the textually-identical "){" belonging
# to the ajax-link anonymous class further down is intentionally left
# unchanged, matching the source diff)
$f = $tr.Find("){", 1)
$f.Text = ") {"

# "if (optional.isPresent()){" -> "... ()) {"
$f = $tr.Find("()){")
$f.Text = "()) {"

# Slide 24, shape 3: "public EchoFormPage(){" -> "public EchoFormPage() {"
$sh = $p.Slides.Item(24).Shapes.Item(3)
$tr = $sh.TextFrame.TextRange
$f = $tr.Find("(){")
$f.Text = "() {"
